$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Product backlog grew by three new items (IDs 10, 11, 12) - add them below
# the existing "Monitoramento do hardware e SO" row.
# ---------------------------------------------------------------------------

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Deverá ter um meio de validação com o cliente"
$ws.Range("C12").Value = "Essencial"

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Deverá conter uma inovação"
$ws.Range("C13").Value = "Essencial"

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "Deverá ser um KPI de acordo com o negócio"
$ws.Range("C14").Value = "Importante "

# ---------------------------------------------------------------------------
# Give the three new rows the same look as the rest of the table: centered
# Arial text inside a thin black box border.
# ---------------------------------------------------------------------------

$newRows = $ws.Range("A12:C14")
$newRows.Font.Name = "Arial"
$newRows.HorizontalAlignment = -4108
$newRows.Borders.LineStyle = 1
$newRows.Borders.Weight = 2

# The table's closing border now sits under row 14, so row 10 no longer
# needs to draw its own bottom edge (row 11 already supplies the divider).
$ws.Range("A10:B10").Borders.Item(9).LineStyle = -4142

# ---------------------------------------------------------------------------
# Leave the cursor where the author finished typing.
# ---------------------------------------------------------------------------

$ws.Range("B16").Select()
